$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.628.41"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.57%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.849.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.65%  "

# Row 4
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.32%  "

# Row 6
$ws.Range("E6").Value = "  +0.03%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4254"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.66%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3645"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.14%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.54"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.89%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07292"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.14%  "

# Row 11
$ws.Range("E11").Value = "  -2.95%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.58"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.07%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.840.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.34%  "

# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.517"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.85%  "

# Row 15
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.313"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.25%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06881"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.30%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.21%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "79.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.11%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008972"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.27%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.11%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.17%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.638.80"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.55%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.985"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.06%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.50%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.081.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.76%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.982"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.25%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.51%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.18%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "121.91"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +10.22%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.263"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.65%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.870"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +12.94%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08863"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.14%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7695"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.65%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.540"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.81%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.966"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.41%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.107"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.55%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05391"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.94%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.094"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.17%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01941"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.75%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.826"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.13%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.872"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.62%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5067"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.31%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1650"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.24%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.374"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.97%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.06530"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.22%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.69%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "104.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.10%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4654"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.20%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.9996"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.02%  "

# Row 50
$ws.Range("E50").Value = "  +0.01%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.38"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.22%  "
